$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 3
    3 = -6
    4 = -3
    5 = 3
    6 = 6
    7 = 1
    8 = -2
    10 = -1
    11 = -4
    12 = 5
    13 = -1
    14 = 2
    15 = 2
    16 = -2
    17 = -3
    18 = 3
    20 = 6
    21 = 4
    23 = -3
    24 = -1
    25 = -1
    26 = -5
    28 = 8
    29 = -2
    30 = 1
    31 = -2
    32 = -1
    33 = 1
    34 = -1
    35 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
